$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# 1. Rename label "Abteilungsabkürzung" -> "Abkuerzung Abteilung" (row 27, col A)
$ws.Range("A27").Value2 = "Abkuerzung Abteilung"

# 2. Insert a new row before row 31 (Gesellschaft), pushing it (and everything below) down by one.
$ws.Rows.Item(31).Insert()

# Copy formatting for the new row 31 (A:B) from the plain "section label" style (A2:B2)
$ws.Range("A2:B2").Copy()
$ws.Range("A31:B31").PasteSpecial(-4122)  # xlPasteFormats

# Copy D-column format (s=10) down into the two Gesellschaft rows (31 & 32)
$ws.Range("D30").Copy()
$ws.Range("D31:D32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Fill content for the (now 2-row) Gesellschaft block
$ws.Range("A31").Value2 = "Gesellschaft"
$ws.Range("B31").Value2 = "Berliner Stadtreinigung"
$ws.Range("A32").Value2 = "Abkuerzung Gesellschaft"
$ws.Range("B32").Value2 = "BSR"

# 4. Update sheet view (scroll position / active selection) to match the saved state
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C19").Select()

$wb.Save()
